# Update the simulated-game transition-probability matrix on Sheet1.
# Each non-zero cell below represents the fraction of simulated games
# that transitioned from the row's starting state to the column's state.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "B2" = 0.4166666666666667
    "C2" = 0.3333333333333333
    "P2" = 0.25
    "P3" = 0.5
    "S3" = 0.5
    "B6" = 0.3333333333333333
    "J6" = 0.1666666666666667
    "S6" = 0.5
    "B7" = 0.2
    "Q7" = 0.2
    "R7" = 0.2
    "S7" = 0.4
    "B8" = 0.07407407407407407
    "F8" = 0.07407407407407407
    "J8" = 0.1481481481481481
    "O8" = 0.03703703703703703
    "Q8" = 0.1481481481481481
    "R8" = 0.03703703703703703
    "S8" = 0.4814814814814815
    "B9" = 0.2
    "J9" = 0.2
    "O9" = 0.2
    "S9" = 0.4
    "B10" = 0.04
    "F10" = 0.06
    "J10" = 0.12
    "Q10" = 0.3
    "R10" = 0.04
    "S10" = 0.44
    "G11" = 0.2
    "J11" = 0.1
    "K11" = 0.2
    "L11" = 0.5
    "G12" = 0.6
    "J12" = 0.2
    "S12" = 0.2
    "H15" = 0.1
    "J15" = 0.6
    "O15" = 0.1
    "S15" = 0.2
    "H16" = 0.4
    "J16" = 0.2
    "O16" = 0.4
    "H17" = 0.25
    "I17" = 0.15
    "J17" = 0.45
    "K17" = 0.1
    "O17" = 0.05
    "H18" = 0.25
    "K18" = 0.25
    "S18" = 0.5
    "H19" = 0.3518518518518519
    "I19" = 0.03703703703703703
    "J19" = 0.3703703703703703
    "K19" = 0.09259259259259259
    "O19" = 0.05555555555555555
    "S19" = 0.09259259259259259
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

